# TablaTFG.xlsx update
# - Renamed task "Implementación API Videojuegos y Precios" -> "Implementación API Videojuegos."
# - Added a new task "Refactorización de código, pequeños ajustes y quitar hardcode" replacing the
#   old "Automatización para actualizar precios" entry at row 15, and re-pointed priorities/times.
# - Row 13 (Exportación de biblioteca) priority bumped Media -> Alta.
# - Rows 14-17 time-tracking values updated.
# - Totals row (18) now sums D13:D17 / E13:E17 instead of D5,D8,D12:D17 / E5,E8,E12:E17.
# - Rows 5, 8 and 12 recolored to match the "in progress" (blue) highlight used by rows 9-11.
# - Column B widened, and the active selection moved to E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Text / priority updates -------------------------------------------------
$ws.Range("B8").Value = "Implementación API Videojuegos."
$ws.Range("C13").Value = "Alta"
$ws.Range("B15").Value = "Refactorización de código, pequeños ajustes y quitar hardcode"
$ws.Range("C15").Value = "Alta"

# --- Recolor rows 5, 8 and 12 to match rows 9-11's fill/format --------------
$ws.Range("A9:C9").Copy()
$ws.Range("A5:C5").PasteSpecial($xlPasteFormats)
$ws.Range("A8:C8").PasteSpecial($xlPasteFormats)
$ws.Range("A12:C12").PasteSpecial($xlPasteFormats)

$ws.Range("D9:E9").Copy()
$ws.Range("D5:E5").PasteSpecial($xlPasteFormats)
$ws.Range("D8:E8").PasteSpecial($xlPasteFormats)
$ws.Range("D12:E12").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Updated tracked time values ---------------------------------------------
$ws.Range("D14").Value = 0.83333333333333337
$ws.Range("E14").Value = 0.083333333333333329

$ws.Range("D15").Value = 0.20833333333333334

$ws.Range("D16").Value = 0.41666666666666669

$ws.Range("D17").Value = 0.41666666666666669

# --- Totals now only cover rows 13-17 ----------------------------------------
$ws.Range("D18").Formula = "=SUM(D13:D17)"
$ws.Range("E18").Formula = "=SUM(E13:E17)"

# --- Layout tweaks -------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 68.16666666666667
$ws.Range("E13").Select()
